# Apply the edits described by the diff:
# 1. Mark several existing "Graph" topics in column C as "Yes" (were "<->").
# 2. Change formatting (smaller hyperlink font + row height) for row 384.
# 3. Insert a new row before row 400 for a new "Graph" topic:
#    "Hamiltonian cycle" (no hyperlink), which shifts every row from the
#    old row 400 onward down by one (Trie / Dynamic Programming / Bit
#    Manipulation sections, and all their hyperlinks, move from rows
#    400-481 to 401-482).
# 4. Update the active sheet view/selection to reflect where the user was
#    working (row 378 visible, C400 selected).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Mark rows as completed ("Yes") in column C -------------------------
$doneRows = @(362, 371, 372, 373, 374, 376, 377, 382, 383, 384)
foreach ($r in $doneRows) {
    $ws.Cells.Item($r, 3).Value = "Yes"
}

# --- 2. Reformat row 384's topic cell (smaller hyperlink font) -------------
$ws.Cells.Item(384, 2).Font.Size = 12
$ws.Rows.Item(384).RowHeight = 19.5

# --- 3. Insert the new "Hamiltonian cycle" row at row 400 ------------------
$ws.Rows.Item(400).Insert()
$ws.Cells.Item(400, 1).Value = "Graph"
$ws.Cells.Item(400, 2).Value = "Hamiltonian cycle"
$ws.Cells.Item(400, 3).Value = "<->"
$ws.Rows.Item(400).RowHeight = 21

# --- 4. Update sheet view / selection ---------------------------------------
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 378
$win.ScrollColumn = 1
$ws.Range("A400:C400").Select()
